$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B5").Font.Size = 11
$ws.Range("B5").Font.Size = 12
Write-Output "done"
